# Update triangle board to put LED centroid at footprint origin
# (adjusts Mid X / Mid Y for LED1-LED3 rows, and moves the saved cursor
# selection, matching the authored OOXML diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LED1 (row 5): Mid X / Mid Y
$ws.Range("B5").Value = 153.331603
$ws.Range("C5").Value = -136.15

# LED2 (row 6): Mid X / Mid Y
$ws.Range("B6").Value = 106.668395
$ws.Range("C6").Value = -136.15

# LED3 (row 7): Mid Y only
$ws.Range("C7").Value = -95.7

# Restore the cursor/selection position saved in the sheet view
$null = $ws.Range("E23").Select()
